# The diff turns:
#   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>
# into four separate runs:
#   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>
#   <w:r><w:t xml:space="preserve"> (</w:t></w:r>
#   <w:r><w:t>Changed main</w:t></w:r>
#   <w:r><w:t>)</w:t></w:r>
#
# A plain Range.InsertAfter() call appends to the existing run whenever the
# inserted text's formatting matches its neighbour, which would merge
# everything back into a single run. Wrapping the three inserts in
# tracked-changes (and accepting each resulting revision individually,
# instead of via AcceptAllRevisions) keeps each insertion as its own run in
# the saved OOXML while leaving the rest of the document untouched.

$d = $word.ActiveDocument
$d.TrackRevisions = $true

$rng = $d.Content
$found = $rng.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'This is a Microsoft word document.'"
}

# Collapse to an insertion point right after the matched text (i.e. before
# the paragraph mark) and append the three new runs in order.
$rng.Collapse(0)
$rng.InsertAfter(" (")

$rng.Collapse(0)
$rng.InsertAfter("Changed main")

$rng.Collapse(0)
$rng.InsertAfter(")")

foreach ($rev in $d.Revisions) {
    $rev.Accept()
}
$d.TrackRevisions = $false
